$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Introduction")

# --- Make room: shift old legend / colour-key rows (old A3:A14) down by 12 --
$ws.Range("A3:A14").EntireRow.Insert()
# old blank rows A9:A14 are now at A21:A26 and are not needed any more
$ws.Range("A21:A26").EntireRow.Delete()

# --- Row 2 (blank spacer under the title): now uses the title style (bold 14, wrap) --
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").Font.Size = 14
$ws.Range("A2").Font.Color = 0
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 18.75

# --- Row 3: big introduction paragraph --
$ws.Range("A3").Value = 'This is the configuration file used to describe various changeable parameters of the process. You should use this file to store settings that are environment related (like paths to programs or resources), user related (email account names, credential names), or plain data (URL of website or name of SAP report to execute). Below, the purpose of each sheet is explained in more detail.'
$ws.Range("A3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 60

# --- Row 4: section header ("Settings") ---
$ws.Range("A4").Value = 'Settings'
$ws.Range("A4").Font.Bold = $true
$ws.Range("A4").Font.Size = 12
$ws.Range("A4").Font.Color = 0
$ws.Range("A4").Font.Name = "Calibri"
$ws.Range("A4").HorizontalAlignment = -4108
$ws.Range("A4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 15.75

# --- Row 5: section description ---
$ws.Range("A5").Value = 'This sheet is the place to store plain data, as well as most user data with the important exception of credential names.'
$ws.Range("A5").Font.Bold = $false
$ws.Range("A5").Font.Size = 11
$ws.Range("A5").Font.Color = 0
$ws.Range("A5").Font.Name = "Calibri"
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").WrapText = $true

# --- Row 6: section header ("Credentials") ---
$ws.Range("A6").Value = 'Credentials'
$ws.Range("A6").Font.Bold = $true
$ws.Range("A6").Font.Size = 12
$ws.Range("A6").Font.Color = 0
$ws.Range("A6").Font.Name = "Calibri"
$ws.Range("A6").HorizontalAlignment = -4108
$ws.Range("A6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 15.75

# --- Row 7: section description ---
$ws.Range("A7").Value = 'The credentials sheet is the place to store your credential names.There is also one special credential, that needs to be defined only once, and which is comprised of the URL, TenancyName and CredentialName required to authenticate to the Orchestrator server using REST API. This is only used when working with QueueItems.'
$ws.Range("A7").Font.Bold = $false
$ws.Range("A7").Font.Size = 11
$ws.Range("A7").Font.Color = 0
$ws.Range("A7").Font.Name = "Calibri"
$ws.Range("A7").HorizontalAlignment = -4131
$ws.Range("A7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 45

# --- Row 8: section header ("Workblocks") ---
$ws.Range("A8").Value = 'Workblocks'
$ws.Range("A8").Font.Bold = $true
$ws.Range("A8").Font.Size = 12
$ws.Range("A8").Font.Color = 0
$ws.Range("A8").Font.Name = "Calibri"
$ws.Range("A8").HorizontalAlignment = -4108
$ws.Range("A8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 15.75

# --- Row 9: section description ---
$ws.Range("A9").Value = 'The workblock names are of the states in the framework. Define the names of workblocks you create here.'
$ws.Range("A9").Font.Bold = $false
$ws.Range("A9").Font.Size = 11
$ws.Range("A9").Font.Color = 0
$ws.Range("A9").Font.Name = "Calibri"
$ws.Range("A9").HorizontalAlignment = -4131
$ws.Range("A9").WrapText = $true

# --- Row 10: section header ("Constants") ---
$ws.Range("A10").Value = 'Constants'
$ws.Range("A10").Font.Bold = $true
$ws.Range("A10").Font.Size = 12
$ws.Range("A10").Font.Color = 0
$ws.Range("A10").Font.Name = "Calibri"
$ws.Range("A10").HorizontalAlignment = -4108
$ws.Range("A10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 15.75

# --- Row 11: section description ---
$ws.Range("A11").Value = 'Typically there is not much for you to add here, although you want to check/change the settings of the Retry mechanism implemented in at the framework layer, during transaction processing, exception recovery, and continuous failiure. Also stores constants used throughout the program, like preconfiguered delays, timeouts. '
$ws.Range("A11").Font.Bold = $false
$ws.Range("A11").Font.Size = 11
$ws.Range("A11").Font.Color = 0
$ws.Range("A11").Font.Name = "Calibri"
$ws.Range("A11").HorizontalAlignment = -4131
$ws.Range("A11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 45

# --- Row 12: section header ("Assets") ---
$ws.Range("A12").Value = 'Assets'
$ws.Range("A12").Font.Bold = $true
$ws.Range("A12").Font.Size = 12
$ws.Range("A12").Font.Color = 0
$ws.Range("A12").Font.Name = "Calibri"
$ws.Range("A12").HorizontalAlignment = -4108
$ws.Range("A12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 15.75

# --- Row 13: section description ---
$ws.Range("A13").Value = 'This sheet is used to fetch assets from Orchestrator. The column name is the key, while the column asset hoolds the asset name in Orchestrator. If there is another local key with the same name, it will be overwritten by the value fetched from Orchestrator.'
$ws.Range("A13").Font.Bold = $false
$ws.Range("A13").Font.Size = 11
$ws.Range("A13").Font.Color = 0
$ws.Range("A13").Font.Name = "Calibri"
$ws.Range("A13").HorizontalAlignment = -4131
$ws.Range("A13").WrapText = $true
$ws.Rows.Item(13).RowHeight = 30

# --- Row 14: trailing blank line, same body style ---
$ws.Range("A14").Font.Bold = $false
$ws.Range("A14").Font.Size = 11
$ws.Range("A14").Font.Color = 0
$ws.Range("A14").Font.Name = "Calibri"
$ws.Range("A14").HorizontalAlignment = -4131
$ws.Range("A14").WrapText = $true

# --- Row 15: legend title text updated (style already correct, bold 11) ---
$ws.Range("A15").Value = '####  Legend of Key Value pair colours####'

# --- Row 16: "you may want to mark keys..." now plain wrap, no horizontal align ---
$ws.Range("A16").Font.Bold = $false
$ws.Range("A16").Font.Size = 11
$ws.Range("A16").Font.Color = 0
$ws.Range("A16").Font.Name = "Calibri"
$ws.Range("A16").WrapText = $true

# --- Column width & selection ---
$ws.Columns.Item(1).ColumnWidth = 117.17
$ws.Range("A3").Select()
